$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.384.38'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').Value = '2.842.13'
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '361.53'
$ws.Range('E5').Value = '  +3.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '112.93'
$ws.Range('E6').Value = '  -2.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  +4.03%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.610'
$ws.Range('E9').Value = '  +3.27%  '
$ws.Range('E10').Value = '  -2.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0872'
$ws.Range('E11').Value = '  +1.37%  '
$ws.Range('E12').Value = '  +1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.15'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.85'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '3.289.92'
$ws.Range('E15').Value = '  +1.64%  '
$ws.Range('D16').Value = '2.866.03'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.936'
$ws.Range('E17').Value = '  +5.00%  '
$ws.Range('D18').Value = '52.288.97'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('E19').Value = '  +3.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.15'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.53'
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('E22').Value = '  +2.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '273.14'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.67'
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('E25').Value = '  +2.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.02'
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.38'
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('E30').Value = '  +2.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0486'
$ws.Range('E31').Value = '  +9.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '35.57'
$ws.Range('E32').Value = '  +3.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '52.50'
$ws.Range('E33').Value = '  +4.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.93'
$ws.Range('E34').Value = '  +2.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.61'
$ws.Range('E35').Value = '  +14.07%  '
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('E38').Value = '  +2.66%  '
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.58'
$ws.Range('E40').Value = '  -0.39%  '
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.55'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '126.99'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.06'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.29'
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.098.97'
$ws.Range('E46').Value = '  +1.93%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.37'
$ws.Range('E47').Value = '  +1.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.29'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('E49').Value = '  +5.87%  '
$ws.Range('E50').Value = '  +1.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.28'
$ws.Range('E51').Value = '  +3.18%  '
